$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet6")
$ws.Range("A4").Value = "l7"
$ws.Range("B4").Value = 0
